$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Acolyte', 'Anthropologist', 'Archaeologist', 'Athlete', 'Azorius Functionary', 'Black Fist Double Agent', 'Boros Legionnaire', 'Caravan Specialist', 'Celebrity Adventurer’s Scion', 'Charlatan', 'City Watch', 'Clan Crafter', 'Cloistered Scholar', 'Cormanthor Refugee', 'Courtier', 'Criminal', 'Dragon Casualty', 'Earthspur Miner', 'Fisher', 'Failed Merchant', 'Faceless', 'Entertainer', 'Far Traveler', 'Folk Hero', 'Gate Urchin', 'Golgari Agent', 'Gruul Anarch', 'Gambler', 'Guild Merchant', 'Guild Artisan', 'Harborfolk', 'Haunted One', 'Hermit', 'Hillsfar Merchant', 'Hillsfar Smuggler', 'House Agent', 'Inheritor', 'Initiate', 'Inquisitor', 'Iron Route Bandit', 'Izzet Engineer', 'Knight', 'Mercenary Veteran', 'Mulmaster Aristocrat', 'Noble', 'Orzhov Representative', 'Outlander', 'Pirate', 'Phlan Insurgent', 'Phlan Refugee', 'Plaintiff', 'Rakdos Cultist', 'Rival Intern', 'Smuggler', 'Spy'),
    @('Insight=Religion', 'Insight=Religion', 'History=Survival', 'Acrobatics=Athletics', 'Insight=Intimidation', 'Deception=Insight', 'Athletics=Intimidation', 'Animal Handling=Survival', 'Perception=Performance', 'Deception=Sleight of Hand', 'Athletics=Insight', 'History=Insight', 'History=Arcana', 'Nature=Survival', 'Insight=Persuasion', 'Deception=Stealth', 'Intimidation=Survival', 'Athletics=Survival', 'History=Survival', 'Investigation=Persuation', 'Deception=Intimidation', 'Acrobatics=Performance', 'Insight=Perception', 'Animal Handling=Survival', 'Deception=Sleight of Hand', 'Nature=Survival', 'Animal Handling=Athletics', 'Deception=Insight', 'Insight=Persuasion', 'Insight=Persuasion', 'Athletics=Sleight of Hand', 'Arcana=Investigation', 'Medicine=Religion', 'Insight=Persuasion', 'Perception=Stealth', 'Investigation=Persuasion', 'Survival=History', 'Athletics=Intimidation', 'Investigation=Religion', 'Stealth=Animal Handling', 'Arcana=Investigation', 'History=Persuasion', 'Athletics=Persuasion', 'Deception=Performance', 'History=Persuasion', 'Intimidation=Religion', 'Survival=Athletics', 'Athletics=Perception', 'Stealth, Survival', 'Insight=Athletics', 'Medicine=Persuasion', 'Acrobatics=Performance', 'History=Investigation', 'Athletics=Deception', 'Deception=Stealth'),
    @('None', 'None', 'Navigator’s tools', 'Vehicles (Land)', 'None', 'Disguise Kit=*Artisan’s Tools', '*Gaming Sets', 'Vehicles (Land)', 'Disguise Kit', 'Disguise Kit=Forgery Kit', 'None', '*Artisan’s Tools', 'None', '*Artisan’s Tools', 'None', '*Gaming Sets=Thieves’ Tools', 'Vehicles (Water)', 'None', 'None', '*Artisan’s Tools', 'Disguise Kit', 'Disguise Kit=*Musical Instruments', '*Musical Instruments', '*Artisan’s Tools=Vehicles (Land)', 'Thieves’ tools=*Musical Instruments', 'Poisoner’s kit', 'Herbalism Kit', '*Gaming Sets', 'Navigator’s Tools', '*Artisan’s Tools', '*Gaming Sets=Vehicles (Water)', 'None', 'Herbalism kit', 'Vehicles (Land)=Vehicles (Water)', 'Forgery kit', 'None', '*Musical Instruments', '*Gaming Sets=Vehicles (Land)', 'Thieves’ Tools=*Artisan’s Tools', '*Gaming Sets=Vehicles (Land)', '*Artisan’s Tools', '*Gaming Sets', '*Gaming Sets=Vehicles (Land)', '*Artisan’s tools=*Musical Instruments', '*Gaming Sets', 'None', '*Musical Instruments', 'Navigator’s Tools=Vehicles (Water)', '*Artisan’s tools=Vehicles (Land)', '*Artisan’s tools', '*Artisan’s tools', '*Musical Instruments', '*Artisan’s tools', 'Vehicles (Water)', '*Gaming Sets=Thieves’ Tools'),
    @('Any=Any', 'Any=Any', 'Any', 'Any', 'Any=Any', 'None', 'Celestial=Draconic', 'Any', 'Any=Any', 'None', 'Any=Any', 'Dwarvish', 'Any=Any', 'Elven', 'Any=Any', 'None', 'Draconic', 'Dwarven=Undercommon', 'Any', 'Any', 'Any', 'None', 'Any', 'None', 'None', 'Elvish', 'Draconic', 'Any', 'None', 'Any', 'None', 'Abyssal', 'Any', 'None', 'Any', 'None', 'Any', 'None', 'None', 'None', 'Vedalken', 'Any', 'None', 'None', 'Any', 'Any=Any', 'Any', 'None', 'None', 'Any', 'Any', 'Abyssal', 'Any', 'None', 'None'),
    @('A Holy Symbol=A Prayer Book=5 Sticks of Incense=Vestments=a Set of Common Clothes=A Belt Pouch which is Containing 15 gp', 'A Leather-Bound Diary=A Bottle of Ink=An Ink Pen=A Set of Traveler’s Clothes=One Trinket of Special Significance=A Pouch Containing 10 gp', 'A Wooden Case Containing a Map to a Ruin=A Bullseye Lantern=A Miner’s Pick=A Set of Traveler’s Clothes=A Shovel=A Two-Person Tent=A Trinket Recovered from a Dig Site=A Pouch Containing 25 gp', 'A Leather Ball=A Lucky Charm=A Set of Traveler’s Clothes=A Pouch Containing 10 gp', 'An Azorius Insignia=A Scroll Containing the Text of a Law Important to you=A Bottle of Blue Ink=A Pen=A Set of Fine Clothes=A Belt Pouch Containing 10 gp (Azorius-minted 1-zino coins)', 'Disguise Kit=Common Clothes=A Tears of Virulence Emblem=A Writ of Free Agency Signed by the Lord Regent=*Artisan’s Tools=A Pouch with 15 gp (payment for services rendered)', 'A Boros Insignia=A Feather from an Angel’s Wing=A Tattered Piece of a Boros Banner (A Souvenir from a Famous Battle)=A Set of Common Clothes=A Belt Pouch Containing 2gp (Boros-minted 1-zino coins)', 'A Whip=A Tent=A Regional Map=A Set of Traveling Clothes=Belt Pouch Containing 10 gp', 'Disguise kit=A Set of Fine Clothes=Belt Pouch Containing 30 gp(gold piece)', 'A Set of Fine Clothes=Disguise Kit=Tools of the Con of your Choice (Ten Stoppered Bottles filled with Colored Liquid, a Set of Weighted Dice, a Deck of Marked Cards, or a Signet Ring of an Imaginary Duke)=Belt Pouch Containing 15 gp', 'A Uniform in the Style of your Unit=Indicative of your Rank=Horn with which to Summon Help=A Set of Manacles=A Pouch Containing 10 gp', '*Artisan’s Tools=A Maker’s Mark Chisel used to mark your Handiwork with the Symbol of the Clan of Crafters you learned your Skill from=A Set of Traveler’s Clothes=A Pouch containing 5 gp=A Gem worth 10 gp', 'The Scholar’s Robes of your Cloister=A Writing Kit (Small Pouch with a Quill, Ink, Folded Parchment, and a Small Penknife)=A Borrowed Book on the Subject of your Current Study=A Pouch containing 10 gp', 'A Two-Person Tent=*Artisan’s Tools=A Holy Symbol=A Set of Traveler’s Clothes=A Belt Pouch containing 5 gp', 'A Set of Fine Clothes=A Pouch containing 5 gp', 'A Crowbar=A Set of Dark Common Clothes including a Hood=A Pouch containing 15 gp', 'A Dagger=Tattered Rags=A Loaf of Moldy Bread=A Small Cast-Off Scale belonging to Vorgansharax – the Maimed Virulence=A Pouch with 5 gp of Various Coins (salvaged during your escape from Phlan)', 'A Shovel=A Block and Tackle=A Climber’s Kit=A Set of Common Clothes=A Belt Pouch containing 5 gp.', 'Fishing tackle=A Net=A Favorite Fishing Lure=A Set of Traveler’s Clothes=A Belt Pouch containing 10 gp', '*Artisan’s Tools=Merchant’s Scale=A Set of Fine Clothes=A Belt Pouch containing 10 gp', 'A Disguise Kit=A Costume=A Pouch containing 10gp', '*Musical Instruments=The Favor of an Admirer (Trinket)=A Costume=A Belt Pouch containing 15 gp.', 'One Set of Traveler’s Clothes=*Musical Instruments=Poorly Wrought Maps from your Homeland that Depict where you are in Faerûn=A Small Piece of Jewelry worth 10 gp in the Style of your Homeland’s Craftsmanship=A Pouch containing 5 gp', '*Artisan’s Tools=A Shovel=An Iron Pot=A Set of Common Clothes=A Pouch containing 10 gp', 'A Battered Alms Box=*Musical Instruments=A Cast-Off Military Jacket=A Set of Common Clothes=A Belt Pouch containing 10 gp', 'A Golgari Insignia=A Poisoner’s Kit=A Pet Beetle=A Set of Common Clothes=A Belt Pouch containing 10 gp worth of Mixed Coins', 'A Gruul Insignia=A Hunting Trap=An Herbalism Kit=The Skull of a Boar=A Beast-Hide Cloak=A Set of Traveler’s Clothes=A Pouch containing 10gp (Azorius 1-zino coins)', '*Gaming Sets=A Lucky Charm=A Set of Fine Clothes=A Belt Pouch containing 15 gp', 'A Mule=A Cart=A Letter of Introduction from your Guild=A Set of Traveler’s Clothes=A Belt Pouch containing 15 gp', '*Artisan’s Tools=A Letter of Introduction from your Guild=A Set of Traveler’s Clothes=A Belt Pouch containing 15 gp', 'Fishing tackle=*Gaming Sets=A Set of Common Clothes=Rowboat=A Belt containing 5 gp', 'Monster Hunter’s Pack=One Trinket of Special Significance', 'A Scroll Case stuffed Full of Notes from your Studies=A Winter Blanket=A Set of Common Clothes=An Herbalism Kit=5 gp', 'A Fine Set of Clothes=A Signet Ring=A Letter of Introduction from your Family’s Trading House=A Purse containing 25 gp', 'A Forgery Kit=A Set of Common Clothes=A Belt Pouch containing 5 gp', 'A Set of Fine Clothes=House Signet Ring=Identification Papers=A Purse containing 20 gp', 'Your Inheritance=A Set of Traveler’s Clothes=*Musical Instruments=A Pouch containing 15 gp', 'A Simple Puzzle Box=A Scroll containing the Basic Teachings of the Five Gods=*Gaming Sets=A Set of Common Clothes=A Pouch containing 15 gp', 'A Holy Symbol=A Set of Traveler’s Clothes=A Pouch containing 15 gp', 'A Set of Dark Common Clothes=Pack Saddle=Burglar’s Pack=A Pouch containing 5 gp', 'An Izzet Insignia=*Artisan’s Tools=The Charred and Twisted remains of a Failed Experiment=A Hammer=A Block and Tackle=A Set of Common Clothes=A Belt Pouch containing 5 gp (Azorius 1-Zino Coins)', 'A Set of Fine Clothes=A Signet Ring=A Scroll of Pedigree=A Purse containing 25 gp', 'A Uniform of your Company (Traveler’s Clothes in Quality)=An Insignia of your Rank=*Gaming Sets=A Pouch containing 10 gp', '*Artisan’s tools=A Set of Fine Clothes=A Purse containing 10 gp', 'A set of Fine Clothes=A Signet Ring=A Scroll of Pedigree=A Purse containing 25 gp', 'An Orzhov Insignia=A Foot-Long Chain made of Ten Gold Coins=Vestments=A Set of Fine Clothes=A Pouch containing 1pp (an Orzhov-minted 10-Zino Coin)', 'A Staff=A Trophy from an Animal you killed=A Belt Pouch containing 10 gp=A Hunting Trap=A Set of Traveler’s Clothes.', 'A Belaying Pin (club)=50 feet of Silk Rope=A Lucky Charm=A Set of Common Clothes=A Belt Pouch containing 10 gp', 'A bag of 20 Caltrops=A Small Trinket that connects you to the Life you once had before the Occupation of Phlan=A Healer’s Kit=A Set of Dark Common Clothes that includes a Cloak and Hood=A Pouch containing 5 gp', '*Artisan’s tools=A Token of the Life you once knew=A Set of Traveler’s Clothes=A Pouch containing 15 gp', '*Artisan’s tools=Fine Clothes=20 gp', 'A Rakdos Insignia=*Musical Instruments=A Costume=A Hooded Lantern made of Wrought Iron=A 10-foot Length of Chain with Sharply Spiked Links=A tinderbox=10 torches=A Set of Common Clothes=A Belt Pouch containing 10 gp (A Mix of Azorius and Boros 1-zino coins)=A Bottle of Sweet, Red Juice', '*Artisan’s tools=A Ledger from your Previous Employer containing a Small Piece of Useful Information=A Set of Fine Clothes=A Belt Pouch containing 10 gp', 'A Fancy Leather Vest=A Set of Common Clothes=A Leather Pouch with 15 gp', 'A crowbar=A Set of Dark Common Clothes including a Hood=A Pouch containing 15 gp')
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

Write-Output "Done"